$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pistures")
$ws.Activate()

# Rename the existing "Bubble" asset in row 9 to "Bubble Blue"
$ws.Range("A9").Value = "Bubble Blue"

# Add a new row 10 for the second bubble asset ("Bubble gray"),
# re-using the same "Bubble" asset-name/type/style as row 9.
$ws.Range("A10").Value = "Bubble gray"
$ws.Range("B10").Value = "Bubble"
$ws.Range("C10").Value = "No Attribution"
$ws.Range("E10").Value = "http://photobucket.com/terms"
$ws.Range("F10").Value = "http://smg.photobucket.com/user/PrayciousAnjel/media/Fantasy%20Crests/Fantasy%20Orbs/glassorb1.png.html"

# Match the "No Attribution" styling used in column C for the other rows
$ws.Range("C9").Copy()
$ws.Range("C10").PasteSpecial(-4122)

# Leave the new cell selected, as in the authored workbook
$ws.Range("B10").Select() | Out-Null
